$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.266.32'
$ws.Range("E2").Value = '  +5.67%  '
$ws.Range("D3").Value = '3.522.10'
$ws.Range("E3").Value = '  +3.33%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'595.98"
$ws.Range("E5").Value = '  +4.90%  '
$ws.Range("D6").Value = "'170.64"
$ws.Range("E6").Value = '  +8.66%  '
$ws.Range("D8").Value = '3.522.04'
$ws.Range("E8").Value = '  +3.21%  '
$ws.Range("E9").Value = '  +1.74%  '
$ws.Range("D10").Value = "'7.28"
$ws.Range("E11").Value = '  +6.21%  '
$ws.Range("D12").Value = "'0.441"
$ws.Range("E12").Value = '  +4.36%  '
$ws.Range("D13").Value = '4.129.02'
$ws.Range("E13").Value = '  +3.29%  '
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").Value = "'28.34"
$ws.Range("E15").Value = '  +5.39%  '
$ws.Range("E16").Value = '  +6.38%  '
$ws.Range("D17").Value = '67.123.25'
$ws.Range("D18").Value = '3.518.21'
$ws.Range("E18").Value = '  +4.27%  '
$ws.Range("D19").Value = "'6.32"
$ws.Range("E19").Value = '  +3.70%  '
$ws.Range("D20").Value = "'14.13"
$ws.Range("E20").Value = '  +4.21%  '
$ws.Range("D21").Value = "'399.06"
$ws.Range("E21").Value = '  +3.85%  '
$ws.Range("D22").Value = "'7.93"
$ws.Range("B23").Value = 'PEPE'
$ws.Range("C23").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D23").Value = "'0.0000130"
$ws.Range("E23").Value = '  +14.68%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = "'73.70"
$ws.Range("E24").Value = '  +3.66%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = "'0.532"
$ws.Range("E26").Value = '  +3.59%  '
$ws.Range("D27").Value = "'10.12"
$ws.Range("E27").Value = '  +4.67%  '
$ws.Range("E28").Value = '  +2.52%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").Value = "'6.41"
$ws.Range("E30").Value = '  +5.95%  '
$ws.Range("D31").Value = "'1.49"
$ws.Range("E31").Value = '  +7.51%  '
$ws.Range("E32").Value = '  +4.57%  '
$ws.Range("D33").Value = "'23.67"
$ws.Range("E33").Value = '  +3.48%  '
$ws.Range("E34").Value = '  +6.92%  '
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("E36").Value = '  +5.78%  '
$ws.Range("D37").Value = "'162.03"
$ws.Range("E37").Value = '  +0.86%  '
$ws.Range("E38").Value = '  +8.27%  '
$ws.Range("E39").Value = '  +7.35%  '
$ws.Range("D40").Value = "'0.0751"
$ws.Range("E40").Value = '  +3.98%  '
$ws.Range("D41").Value = "'4.70"
$ws.Range("E41").Value = '  +8.15%  '
$ws.Range("E42").Value = '  +5.57%  '
$ws.Range("D43").Value = "'26.50"
$ws.Range("E43").Value = '  +2.34%  '
$ws.Range("D44").Value = '2.830.17'
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("D45").Value = "'26.94"
$ws.Range("E45").Value = '  +5.24%  '
$ws.Range("D46").Value = "'43.60"
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("D47").Value = "'2.58"
$ws.Range("E47").Value = '  +10.55%  '
$ws.Range("E48").Value = '  +4.51%  '
$ws.Range("D49").Value = "'353.31"
$ws.Range("E49").Value = '  +7.74%  '
$ws.Range("E50").Value = '  +6.87%  '
$ws.Range("D51").Value = "'33.50"
$ws.Range("E51").Value = '  +11.54%  '
